$d = $word.ActiveDocument
$nbsp = [char]0xA0

# --- 1) Collapse the three-run templated paragraphs into single runs -------
# "Nombre de malades total : ..." paragraph (was split across 3 <w:r> runs)
$t1 = "Nombre de malades total" + $nbsp + ": {% if object.nb_sick_persons is not none %}{{ object.nb_sick_persons }}{% else %}-{% endif %}"
$d.Content.Find.Execute($t1, $true, $false, $false, $false, $false, $true, 1, $false, $t1, 2)

# "Dont conduits a l'hopital : ..." paragraph
$t2 = "Dont conduits " + [char]0xE0 + " l" + [char]0x2019 + "hopital" + $nbsp + ": {% if object.nb_sick_persons_to_hospital is not none %}{{ object.nb_sick_persons_to_hospital }}{% else %}-{% endif %}"
$d.Content.Find.Execute($t2, $true, $false, $false, $false, $false, $true, 1, $false, $t2, 2)

# "Dont decedes : ..." paragraph
$t3 = "Dont d" + [char]0xE9 + "c" + [char]0xE9 + "d" + [char]0xE9 + "s" + $nbsp + ": {% if object.nb_dead_persons is not none %}{{ object.nb_dead_persons }}{% else %}-{% endif %}"
$d.Content.Find.Execute($t3, $true, $false, $false, $false, $false, $true, 1, $false, $t3, 2)

# --- 2) Insert "Autre identifiant" paragraph right after "N SIRET" --------
$siretText = "N" + [char]0xB0 + " SIRET" + $nbsp + ": {{ etablissement.siret or '-'  }}"

$siretIdx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text -eq ($siretText + [char]13)) {
        $siretIdx = $i
        break
    }
}

$d.Paragraphs.Item($siretIdx).Range.InsertParagraphAfter()

$newText = "Autre identifiant" + $nbsp + ": {{ etablissement.autre_identifiant or '-'  }}"
$newIdx = $siretIdx + 1
$newPara = $d.Paragraphs.Item($newIdx)
$newPara.Range.Text = $newText
